$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts Late/heading/Outstanding
# columns one position to the right: N->O, O->P, P->Q). Excel's default
# column insert copies formatting (incl. width) from the column to the left,
# so mirror column M's width onto the freshly inserted column N.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Update the selection / view on the Repayment schedule sheet and make it
# the active tab of the workbook.
$ws.Range("S8").Select()
$ws.Activate()
